# Updates data in project list template
#
# The commit changes the "Non-Erection Wind Delay Critical Height (m)"
# value in cell AD2 from 12.2 to 11, and leaves the sheet scrolled/
# selected around that same cell (selection moves to AD2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Data edit: AD2 12.2 -> 11
$ws.Range("AD2").Value = 11

# Reflect the author's final cursor position / view state around the
# edited cell (matches the workbook's saved selection: AD2).
$ws.Range("AD2").Select()
